$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 43.995596
$ws.Range("H2").Value = 131.986788
$ws.Range("I2").Value = 0.08241811124115486
$ws.Range("J2").Value = 0.08241811124115485
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2520896666666667
$ws.Range("N2").Value = 0.7562690000000001
$ws.Range("O2").Value = 0.03491140780587004
$ws.Range("P2").Value = 0.03491140780587004
$ws.Range("Q2").Value = 11.09083513044133
$ws.Range("R2").Value = 99.81751617397201
$ws.Range("S2").Value = 0.002877332292129519
$ws.Range("T2").Value = 0.002877332292129518
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 43.995596
$ws.Range("H3").Value = 131.986788
$ws.Range("I3").Value = 0.08241811124115486
$ws.Range("J3").Value = 0.08241811124115485
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.312792666666667
$ws.Range("N3").Value = 3.938378
$ws.Range("O3").Value = 0.1818061039810792
$ws.Range("P3").Value = 0.1818061039810792
$ws.Range("Q3").Value = 57.75709579442933
$ws.Range("R3").Value = 519.813862149864
$ws.Range("S3").Value = 0.01498411570223356
$ws.Range("T3").Value = 0.01498411570223355
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 43.995596
$ws.Range("H4").Value = 131.986788
$ws.Range("I4").Value = 0.08241811124115486
$ws.Range("J4").Value = 0.08241811124115485
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.655957
$ws.Range("N4").Value = 16.967871
$ws.Range("O4").Value = 0.7832824882130508
$ws.Range("P4").Value = 0.7832824882130508
$ws.Range("Q4").Value = 248.837199165372
$ws.Range("R4").Value = 2239.534792488348
$ws.Range("S4").Value = 0.0645566632467918
$ws.Range("T4").Value = 0.06455666324679178
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 439.8208616666666
$ws.Range("H5").Value = 1319.462585
$ws.Range("I5").Value = 0.8239280291378236
$ws.Range("J5").Value = 0.8239280291378236
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2520896666666667
$ws.Range("N5").Value = 0.7562690000000001
$ws.Range("O5").Value = 0.03491140780587004
$ws.Range("P5").Value = 0.03491140780587004
$ws.Range("Q5").Value = 110.8742944105961
$ws.Range("R5").Value = 997.8686496953651
$ws.Range("S5").Value = 0.02876448742791733
$ws.Range("T5").Value = 0.02876448742791733
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 439.8208616666666
$ws.Range("H6").Value = 1319.462585
$ws.Range("I6").Value = 0.8239280291378236
$ws.Range("J6").Value = 0.8239280291378236
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.312792666666667
$ws.Range("N6").Value = 3.938378
$ws.Range("O6").Value = 0.1818061039810792
$ws.Range("P6").Value = 0.1818061039810792
$ws.Range("Q6").Value = 577.3936018430145
$ws.Range("R6").Value = 5196.54241658713
$ws.Range("S6").Value = 0.1497951449383568
$ws.Range("T6").Value = 0.1497951449383568
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 439.8208616666666
$ws.Range("H7").Value = 1319.462585
$ws.Range("I7").Value = 0.8239280291378236
$ws.Range("J7").Value = 0.8239280291378236
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.655957
$ws.Range("N7").Value = 16.967871
$ws.Range("O7").Value = 0.7832824882130508
$ws.Range("P7").Value = 0.7832824882130508
$ws.Range("Q7").Value = 2487.607881289615
$ws.Range("R7").Value = 22388.47093160653
$ws.Range("S7").Value = 0.6453683967715494
$ws.Range("T7").Value = 0.6453683967715494
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 49.99334866666667
$ws.Range("H8").Value = 149.980046
$ws.Range("I8").Value = 0.09365385962102149
$ws.Range("J8").Value = 0.09365385962102149
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2520896666666667
$ws.Range("N8").Value = 0.7562690000000001
$ws.Range("O8").Value = 0.03491140780587004
$ws.Range("P8").Value = 0.03491140780587004
$ws.Range("Q8").Value = 12.60280660093045
$ws.Range("R8").Value = 113.425259408374
$ws.Range("S8").Value = 0.003269588085823186
$ws.Range("T8").Value = 0.003269588085823186
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 49.99334866666667
$ws.Range("H9").Value = 149.980046
$ws.Range("I9").Value = 0.09365385962102149
$ws.Range("J9").Value = 0.09365385962102149
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.312792666666667
$ws.Range("N9").Value = 3.938378
$ws.Range("O9").Value = 0.1818061039810792
$ws.Range("P9").Value = 0.1818061039810792
$ws.Range("Q9").Value = 65.63090151170978
$ws.Range("R9").Value = 590.678113605388
$ws.Range("S9").Value = 0.01702684334048883
$ws.Range("T9").Value = 0.01702684334048883
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 49.99334866666667
$ws.Range("H10").Value = 149.980046
$ws.Range("I10").Value = 0.09365385962102149
$ws.Range("J10").Value = 0.09365385962102149
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.655957
$ws.Range("N10").Value = 16.967871
$ws.Range("O10").Value = 0.7832824882130508
$ws.Range("P10").Value = 0.7832824882130508
$ws.Range("Q10").Value = 282.760230344674
$ws.Range("R10").Value = 2544.842073102066
$ws.Range("S10").Value = 0.07335742819470947
$ws.Range("T10").Value = 0.07335742819470947
